# ---------------------------------------------------------------
# Add the new "2022-Q3" quarterly holdings sheet, positioned right
# after the "总计" summary sheet, and prepend its summary row.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2022-Q2")

# --- create the new quarter sheet (data/format written in place,
# THEN moved - the sheet reference goes stale once Move() runs) --
$q3 = $wb.Worksheets.Add()
$q3.Name = "2022-Q3"

# copy the header-row / index-column formatting from an existing
# quarterly sheet so the new sheet matches the others
$template.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q3.Range("A2:A16").PasteSpecial(-4122)

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $q3.Cells.Item(1, 2 + $c).Value = $headers[$c]
}

# columns B:G hold text-formatted numbers (matches the source data)
$q3.Range("B2:G16").NumberFormat = "@"

$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "011866"
$q3.Cells.Item(2,3).Value = "广发价值增长混合A"
$q3.Cells.Item(2,4).Value = "15.73"
$q3.Cells.Item(2,5).Value = "94.64"
$q3.Cells.Item(2,6).Value = "6.11"
$q3.Cells.Item(2,7).Value = "0.9611"
$q3.Cells.Item(2,8).Value = 9
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "002624"
$q3.Cells.Item(3,3).Value = "广发优企精选灵活配置混合A"
$q3.Cells.Item(3,4).Value = "11.48"
$q3.Cells.Item(3,5).Value = "94.33"
$q3.Cells.Item(3,6).Value = "7.90"
$q3.Cells.Item(3,7).Value = "0.9069"
$q3.Cells.Item(3,8).Value = 5
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).Value = "270025"
$q3.Cells.Item(4,3).Value = "广发行业领先混合A"
$q3.Cells.Item(4,4).Value = "9.95"
$q3.Cells.Item(4,5).Value = "94.64"
$q3.Cells.Item(4,6).Value = "6.37"
$q3.Cells.Item(4,7).Value = "0.6338"
$q3.Cells.Item(4,8).Value = 8
$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).Value = "501070"
$q3.Cells.Item(5,3).Value = "广发睿阳三年定期开放混合"
$q3.Cells.Item(5,4).Value = "6.62"
$q3.Cells.Item(5,5).Value = "51.01"
$q3.Cells.Item(5,6).Value = "6.48"
$q3.Cells.Item(5,7).Value = "0.4290"
$q3.Cells.Item(5,8).Value = 1
$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).Value = "011427"
$q3.Cells.Item(6,3).Value = "广发价值驱动混合A"
$q3.Cells.Item(6,4).Value = "2.64"
$q3.Cells.Item(6,5).Value = "94.69"
$q3.Cells.Item(6,6).Value = "6.27"
$q3.Cells.Item(6,7).Value = "0.1655"
$q3.Cells.Item(6,8).Value = 8
$q3.Cells.Item(7,1).Value = 5
$q3.Cells.Item(7,2).Value = "200010"
$q3.Cells.Item(7,3).Value = "长城双动力混合A"
$q3.Cells.Item(7,4).Value = "3.29"
$q3.Cells.Item(7,5).Value = "93.10"
$q3.Cells.Item(7,6).Value = "2.97"
$q3.Cells.Item(7,7).Value = "0.0977"
$q3.Cells.Item(7,8).Value = 8
$q3.Cells.Item(8,1).Value = 6
$q3.Cells.Item(8,2).Value = "015561"
$q3.Cells.Item(8,3).Value = "长城双动力混合C"
$q3.Cells.Item(8,4).Value = "2.72"
$q3.Cells.Item(8,5).Value = "93.10"
$q3.Cells.Item(8,6).Value = "2.97"
$q3.Cells.Item(8,7).Value = "0.0808"
$q3.Cells.Item(8,8).Value = 8
$q3.Cells.Item(9,1).Value = 7
$q3.Cells.Item(9,2).Value = "000747"
$q3.Cells.Item(9,3).Value = "广发逆向策略灵活配置混合A"
$q3.Cells.Item(9,4).Value = "1.00"
$q3.Cells.Item(9,5).Value = "94.25"
$q3.Cells.Item(9,6).Value = "5.18"
$q3.Cells.Item(9,7).Value = "0.0518"
$q3.Cells.Item(9,8).Value = 10
$q3.Cells.Item(10,1).Value = 8
$q3.Cells.Item(10,2).Value = "210002"
$q3.Cells.Item(10,3).Value = "金鹰红利价值混合A"
$q3.Cells.Item(10,4).Value = "1.19"
$q3.Cells.Item(10,5).Value = "61.62"
$q3.Cells.Item(10,6).Value = "4.08"
$q3.Cells.Item(10,7).Value = "0.0486"
$q3.Cells.Item(10,8).Value = 5
$q3.Cells.Item(11,1).Value = 9
$q3.Cells.Item(11,2).Value = "011867"
$q3.Cells.Item(11,3).Value = "广发价值增长混合C"
$q3.Cells.Item(11,4).Value = "0.72"
$q3.Cells.Item(11,5).Value = "94.64"
$q3.Cells.Item(11,6).Value = "6.11"
$q3.Cells.Item(11,7).Value = "0.0440"
$q3.Cells.Item(11,8).Value = 9
$q3.Cells.Item(12,1).Value = 10
$q3.Cells.Item(12,2).Value = "011428"
$q3.Cells.Item(12,3).Value = "广发价值驱动混合C"
$q3.Cells.Item(12,4).Value = "0.42"
$q3.Cells.Item(12,5).Value = "94.69"
$q3.Cells.Item(12,6).Value = "6.27"
$q3.Cells.Item(12,7).Value = "0.0263"
$q3.Cells.Item(12,8).Value = 8
$q3.Cells.Item(13,1).Value = 11
$q3.Cells.Item(13,2).Value = "016563"
$q3.Cells.Item(13,3).Value = "金鹰红利价值混合C"
$q3.Cells.Item(13,4).Value = "0.34"
$q3.Cells.Item(13,5).Value = "61.62"
$q3.Cells.Item(13,6).Value = "4.08"
$q3.Cells.Item(13,7).Value = "0.0139"
$q3.Cells.Item(13,8).Value = 5
$q3.Cells.Item(14,1).Value = 12
$q3.Cells.Item(14,2).Value = "010021"
$q3.Cells.Item(14,3).Value = "广发优企精选灵活配置混合C"
$q3.Cells.Item(14,4).Value = "0.07"
$q3.Cells.Item(14,5).Value = "94.33"
$q3.Cells.Item(14,6).Value = "7.90"
$q3.Cells.Item(14,7).Value = "0.0055"
$q3.Cells.Item(14,8).Value = 5
$q3.Cells.Item(15,1).Value = 13
$q3.Cells.Item(15,2).Value = "011758"
$q3.Cells.Item(15,3).Value = "广发逆向策略灵活配置混合C"
$q3.Cells.Item(15,4).Value = "0.02"
$q3.Cells.Item(15,5).Value = "94.25"
$q3.Cells.Item(15,6).Value = "5.18"
$q3.Cells.Item(15,7).Value = "0.0010"
$q3.Cells.Item(15,8).Value = 10
$q3.Cells.Item(16,1).Value = 14
$q3.Cells.Item(16,2).Value = "960001"
$q3.Cells.Item(16,3).Value = "广发行业领先混合H"
$q3.Cells.Item(16,4).Value = "0.01"
$q3.Cells.Item(16,5).Value = "94.64"
$q3.Cells.Item(16,6).Value = "6.37"
$q3.Cells.Item(16,7).Value = "0.0006"
$q3.Cells.Item(16,8).Value = 8

# now that the sheet is fully populated, move it into position
# (right after "总计", i.e. before the old second sheet)
$q3.Move($template)

# ---------------------------------------------------------------
# Prepend the 2022-Q3 row to the "总计" summary sheet, shifting the
# existing quarters down by one row.
# ---------------------------------------------------------------
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 15
$summary.Cells.Item(2,4).Value = 3.47
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q2"
$summary.Cells.Item(3,3).Value = 13
$summary.Cells.Item(3,4).Value = 3.6
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2022-Q1"
$summary.Cells.Item(4,3).Value = 10
$summary.Cells.Item(4,4).Value = 3.1
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q4"
$summary.Cells.Item(5,3).Value = 9
$summary.Cells.Item(5,4).Value = 4.9
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q3"
$summary.Cells.Item(6,3).Value = 5
$summary.Cells.Item(6,4).Value = 1.14
$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2021-Q2"
$summary.Cells.Item(7,3).Value = 4
$summary.Cells.Item(7,4).Value = 2.13
$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(8,2).Value = "2021-Q1"
$summary.Cells.Item(8,3).Value = 3
$summary.Cells.Item(8,4).Value = 0.78
$summary.Cells.Item(9,1).Value = 7
$summary.Cells.Item(9,2).Value = "2020-Q4"
$summary.Cells.Item(9,3).Value = 1
$summary.Cells.Item(9,4).Value = 0.66

